# feat: add visual cues to promo code entries
#
# Inserts two new blocks of localization rows into the tsv_UI_Defaults sheet:
#   1. Four "UI promo code cue" rows (currency / stat pickup toasts) right
#      before the existing "UI wave" row (old row 25).
#   2. Two "UI player stats" rows (Full heal / Health up) right after the
#      existing "tmp psDodgeDelay" row (old row 32, now row 36 after block 1
#      has been inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$IA = [char]0x00A1  # inverted exclamation mark, used by several Spanish strings

# ---------------------------------------------------------------------------
# Block 1: insert 4 rows above the old row 25 ("UI wave") for the new
# "UI promo code cue" entries.
# ---------------------------------------------------------------------------
$ws.Rows("25:28").Insert()

# Bring over the established "data row" look (thin border box, centered
# vertical alignment, wrapped text) by pasting formats from the row just
# above (still a plain, unmerged data row) instead of re-deriving the
# border/fill from scratch - this keeps the shared style table compact,
# exactly like Excel does when you insert rows in the middle of a table.
$ws.Range("A24:G24").Copy()
$ws.Range("A25:G28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$cue = "UI promo code cue"

# Row 25: small currency obtained
$ws.Range("A25").Value = $cue
$ws.Range("B25").Value = "tmp promo small currency"
$ws.Range("C25").Value = "Small currency obtained!"
$ws.Range("D25").Value = "Petite monnaie obtenue!"
$ws.Range("E25").Value = "${IA}Pequeña moneda obtenida!"
$ws.Range("F25").Value = "小さな通貨を手に入れた！"
$ws.Range("G25").Value = "小货币获得！"
$ws.Rows.Item(25).RowHeight = 28.8

# Row 26: big currency obtained
$ws.Range("A26").Value = $cue
$ws.Range("B26").Value = "tmp promo big currency"
$ws.Range("C26").Value = "Big currency otained!"
$ws.Range("D26").Value = "Grosse monnaie obtenue!"
$ws.Range("E26").Value = "Gran moneda obtenida!"
$ws.Range("F26").Value = "大きな通貨を手に入れた！"
$ws.Range("G26").Value = "获得大额货币！"
$ws.Rows.Item(26).RowHeight = 43.2

# Row 27: permanent statistic obtained
$ws.Range("A27").Value = $cue
$ws.Range("B27").Value = "tmp promo perma stat"
$ws.Range("C27").Value = "Permanent statistic obtained!"
$ws.Range("D27").Value = "Statistique permanente obtenue!"
$ws.Range("E27").Value = "${IA}Estadística permanente obtenida!"
$ws.Range("F27").Value = "永久統計取得！"
$ws.Range("G27").Value = "获得永久统计！"
$ws.Rows.Item(27).RowHeight = 43.2

# Row 28: temporary statistic obtained
$ws.Range("A28").Value = $cue
$ws.Range("B28").Value = "tmp promo temp stat"
$ws.Range("C28").Value = "Temporary statistic obtained!"
$ws.Range("D28").Value = "Statistique temporaire obtenue!"
$ws.Range("E28").Value = "${IA}Estadística temporal obtenida!"
$ws.Range("F28").Value = "暫定統計取得！"
$ws.Range("G28").Value = "临时统计得到！"
$ws.Rows.Item(28).RowHeight = 43.2

# C28 keeps the plain (border-less) look from the source edit.
$ws.Range("C28").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# Block 2: insert 2 rows after the row that now holds "tmp psDodgeDelay"
# (old row 32, shifted to row 36 by block 1) for the new player-stat
# entries ("Full heal" / "Health up").
# ---------------------------------------------------------------------------
$ws.Rows("37:38").Insert()

$ws.Range("A36:G36").Copy()
$ws.Range("A37:G38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$stats = "UI player stats"

# Row 37: Full heal
$ws.Range("A37").Value = $stats
$ws.Range("B37").Value = "tmp psFullHeal"
$ws.Range("C37").Value = "Full heal"
$ws.Range("D37").Value = "Soin complet"
$ws.Range("E37").Value = "Totalmente curado"
$ws.Range("F37").Value = "フルヒール"
$ws.Range("G37").Value = "完全治愈"
$ws.Rows.Item(37).RowHeight = 28.8

# Row 38: Health up
$ws.Range("A38").Value = $stats
$ws.Range("B38").Value = "tmp psHealthUp"
$ws.Range("C38").Value = "Health up"
$ws.Range("D38").Value = "Points de vie additionnels"
$ws.Range("E38").Value = "Salud arriba"
$ws.Range("F38").Value = "ヘルスアップ"
$ws.Range("G38").Value = "健康起来"
$ws.Rows.Item(38).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Refresh the view state to roughly match the author's final camera
# position/zoom over the newly expanded table.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("I38").Select()

Write-Output "Added promo-code-cue rows (25-28) and player-stat rows (37-38)"
